# Add season-record columns (Wins / Losses / Ties) to the STL_2013 sheet.
# The original scraper only pulled team statistics, not the season record,
# so this fills in the missing "Wins", "Losses", "Ties" columns (AD, AE, AF)
# for the header row and every data row (2-46) with the team's record
# (97 wins, 65 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting (bold, border, centered) already used by the other
# header cells by copying the format from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the record for every player row (2 through 46)
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
